# Apply "new constraints & data adjustments" edits to the centres workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1 height tweak (header row) ---
$ws.Rows.Item(1).RowHeight = 19.5

# --- Fix unknown "Jour Distribution" value for "Salies du Salat/Mane" (row 25) ---
$ws.Range("D25").Value = "Mercredi"

# --- Tonnage Surgelé (kg) column H adjustments: new per-centre caps/values ---
$ws.Range("H3").Value = 55
$ws.Range("H7").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("H11").Value = 30
$ws.Range("H14").Value = 40
$ws.Range("H15").Value = 30
$ws.Range("H16").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("H20").Value = 30
$ws.Range("H24").Value = 30
$ws.Range("H29").Value = 30
$ws.Range("H30").Value = 30
